# Cotações atualizadas - 2025-12-06
# Append a new row (92) with the latest quotation values, mirroring the
# layout of the existing data rows (date in column A as a number styled
# like the other date cells, and text-valued quotes in columns B:E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 92

# Column A: serial date number, same style as the row above (A91).
$ws.Cells.Item($newRow, 1).Value = 45997
$ws.Cells.Item($newRow, 1).Style = $ws.Cells.Item($newRow - 1, 1).Style
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat

# Columns B:E: textual quote values (comma decimal separator), stored as text.
$ws.Cells.Item($newRow, 2).Value = "22,2192"
$ws.Cells.Item($newRow, 3).Value = "15,9549"
$ws.Cells.Item($newRow, 4).Value = "15,6518"
$ws.Cells.Item($newRow, 5).Value = "15,6518"
